$wb = $excel.ActiveWorkbook

# Add the new worksheet "tc002" right after the existing "tc001" sheet
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item("tc001"))
$newSheet.Name = "tc002"

# Populate data
$newSheet.Range("A1").Value = "projectName"
$newSheet.Range("A2").Value = "SET- DRV"

# Column width (~23.625 chars as authored; engine snaps to its nearest
# internal grid step, 22.8 is the input that lands closest to that target)
$newSheet.Columns.Item(1).ColumnWidth = 22.8

# Selection on the new sheet
$newSheet.Range("A7").Select()

# Make the new sheet the active tab
$newSheet.Activate()
